# RGF_input_file.xlsx edit script
# - Add GPU max-matrix-size input (row 2 C/D, new row 3 "CPU max matrix")
# - Change Lattice value from "BLG" to "MLG"
# - Add "Plot band structure" enable flag (new row)
# - Add a new ribbon-shape data row, tweak existing numbers
# - Repurpose old "BLG" shared string slot -> "Plot band structure" and add
#   three new shared strings ("Define GPU enable and max calculated matrix
#   size", "MLG", "CPU max matrix") -- handled automatically by the engine
#   as new string values are written to cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Insert the two new rows first (this shifts everything below down,
#    and Excel auto-adjusts the MAX(...)/SUM(...) formula ranges and the
#    merged-cell references for us).
# ---------------------------------------------------------------------
$ws.Rows("3:3").Insert()    # new row for "CPU max matrix"
$ws.Rows("10:10").Insert()  # new row for "Plot band structure"

# ---------------------------------------------------------------------
# 2) Row 2 ("Using GPU"): give B2 the Input style, add a max-matrix-size
#    value in C2 and a comment label in D2. Grab the "Input" style from
#    row 4 (Material/Graphene), which keeps that style untouched by the
#    inserts above.
# ---------------------------------------------------------------------
$ws.Range("B4").Copy()
$ws.Range("B2").PasteSpecial(-4122)   # xlPasteFormats - reuse the Input style
$ws.Range("C2").PasteSpecial(-4122)
$ws.Range("B2").Value = $false
$ws.Range("C2").Value = 4000
$ws.Range("D2").Value = "Define GPU enable and max calculated matrix size"

# ---------------------------------------------------------------------
# 3) New row 3: "CPU max matrix" = 1000 (Input style)
# ---------------------------------------------------------------------
$ws.Range("B4").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("A3").Value = "CPU max matrix"
$ws.Range("B3").Value = 1000

# ---------------------------------------------------------------------
# 4) Row 5 ("Lattice"): change value from "BLG" to "MLG"
# ---------------------------------------------------------------------
$ws.Range("B5").Value = "MLG"

# ---------------------------------------------------------------------
# 5) New row 10: "Plot band structure" = TRUE (Input style)
# ---------------------------------------------------------------------
$ws.Range("B4").Copy()
$ws.Range("B10").PasteSpecial(-4122)
$ws.Range("A10").Value = "Plot band structure"
$ws.Range("B10").Value = $true

# ---------------------------------------------------------------------
# 6) Row 14 ("o" data row): F14 300 -> 10, I14 2.5 -> 0
# ---------------------------------------------------------------------
$ws.Range("F14").Value = 10
$ws.Range("I14").Value = 0

# ---------------------------------------------------------------------
# 7) Row 15 ("x" data row, existing): E15 2 -> 1, F15 300 -> 0
# ---------------------------------------------------------------------
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0

# ---------------------------------------------------------------------
# 8) New row 16: duplicate of row 15 ("x" data row), F16 = 0
# ---------------------------------------------------------------------
$ws.Range("A16").Value = "x"
$ws.Range("B16").Value = 1
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0

# ---------------------------------------------------------------------
# 9) Update selection to match the authored file
# ---------------------------------------------------------------------
$ws.Range("F15").Select()
